$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume cells to remain plain text (matching the
# original inline-string cells) instead of being auto-converted to
# numbers by Excel, then restore the default "Normal" style so no
# extra formatting is left behind on the cells.
$targetCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","E25","D26","E26","D27","E27","D28","E28","D29","E29","D30","E30","D31","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","E37","D38","E38","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.056.89"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.839.55"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "245.37"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "0.6969"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").Value = "0.9990"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.07714"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").Value = "23.50"
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").Value = "0.07825"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "92.94"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.841.38"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").Value = "5.116"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "0.6839"
$ws.Range("D16").Value = "6.591"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "0.000008280"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").Value = "29.018.56"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "242.15"
$ws.Range("E19").Value = "  -3.36%  "
$ws.Range("D20").Value = "2.076.74"
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").Value = "12.73"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").Value = "0.9986"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "7.491"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("D26").Value = "158.93"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").Value = "8.797"
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("D28").Value = "18.23"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").Value = "1.542"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("D30").Value = "4.221"
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("D31").Value = "4.165"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").Value = "1.195"
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").Value = "0.05114"
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").Value = "0.7794"
$ws.Range("E34").Value = "  +2.78%  "
$ws.Range("D35").Value = "1.858"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("D36").Value = "1.147"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").Value = "1.297.38"
$ws.Range("E38").Value = "  +6.10%  "
$ws.Range("D39").Value = "0.01863"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "2.703"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D41").Value = "0.9481"
$ws.Range("E41").Value = "  +5.27%  "
$ws.Range("D42").Value = "6.147"
$ws.Range("E42").Value = "  +5.69%  "
$ws.Range("D43").Value = "107.81"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("D44").Value = "0.9985"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "9.684"
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").Value = "0.5173"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").Value = "1.977.56"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").Value = "64.17"
$ws.Range("E48").Value = "  -4.95%  "
$ws.Range("D49").Value = "1.759"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "0.00000000120"
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("D51").Value = "6.987"
$ws.Range("E51").Value = "  -0.73%  "

foreach ($addr in $targetCells) {
    $ws.Range($addr).Style = "Normal"
}
